$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ClassName"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Expected "
$ws.Range("E1").Value = "Actual"

$ws.Range("A2").Value = "Login_Page"
$ws.Range("B2").Value = "Admin"
$ws.Range("C2").Value = "admin123"
$ws.Range("D2").Value = "Valid"
$ws.Range("E2").Value = "Valid"

$ws.Range("A3").Value = "Login_Page"
$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "admin123"
$ws.Range("D3").Value = "Valid"
$ws.Range("E3").Value = "Valid"

$ws.Range("A4").Value = "Login_Page"
$ws.Range("B4").Value = "Password"
$ws.Range("C4").Value = "admin"
$ws.Range("D4").Value = "Invalid"
$ws.Range("E4").Value = "Invalid"

$ws.Range("A5").Value = "Validate_Login"
$ws.Range("B5").Value = "UserName"
$ws.Range("C5").Value = "admin"
$ws.Range("D5").Value = "Invalid"
$ws.Range("E5").Value = "Invalid"

$ws.Range("E5").Select()
